# Generate Report for Handback
# Applies the handback-report update described by the commit:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - New "Latest Target File" (F) and "Latest Handback File" (G) columns are
#    populated (with hyperlinks) on the zh-cn and de-de sheets
#  - "Latest Handback DateTime" (H) is populated with real timestamps
#    instead of the zero-date placeholder

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (B and C) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# Helper to find a hyperlink on a worksheet by its cell address (e.g. "$A$3")
function Get-HyperlinkByAddr($ws, $addrTarget) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addrTarget) {
            return $hl
        }
    }
    return $null
}

# Helper that adds a hyperlink cell and paints it with the same
# underline + font colour used by the existing hyperlink cells
# (single underline, RGB(100,149,237) == 0xED9564 as an OLE colour).
function Add-StyledHyperlink($ws, $cellAddr, $address, $display) {
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $address, "", "", $display)
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = 15570276
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("H2").Value = "2016-03-30 11:28:17"
$wsZh.Range("H3").Value = "2016-03-30 11:28:17"

# Pull A3/D3 hyperlink info so the row-3 handoff links can be re-inserted
# after the new row-2 F/G links -- this reproduces the A2,D2,F2,G2,A3,D3,F3,G3
# ordering that Excel itself produced for this edit.
$hlA3 = Get-HyperlinkByAddr $wsZh '$A$3'
$a3Address = $hlA3.Address
$a3Display = $hlA3.TextToDisplay
$hlA3.Delete()

$hlD3 = Get-HyperlinkByAddr $wsZh '$D$3'
$d3Address = $hlD3.Address
$d3Display = $hlD3.TextToDisplay
$hlD3.Delete()

Add-StyledHyperlink $wsZh "F2" "https://github.com/OpenLocalizationTest/oltest/blob/915c598678366622ed8a59ae44ff321bc4968d68/e2e/27908121-ab1a-42e7-a0fa-0f05bdfe0441.md" "27908121-ab1a-42e7-a0fa-0f05bdfe0441.md"
Add-StyledHyperlink $wsZh "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ecd3281a4de4970c00f958a5185cbd454ee743f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/high/27908121-ab1a-42e7-a0fa-0f05bdfe0441.7677fbf96608daf076d6d3278a6411e68340f762.zh-cn.xlf" "27908121-ab1a-42e7-a0fa-0f05bdfe0441.7677fbf96608daf076d6d3278a6411e68340f762.zh-cn.xlf"

Add-StyledHyperlink $wsZh "A3" $a3Address $a3Display
Add-StyledHyperlink $wsZh "D3" $d3Address $d3Display

Add-StyledHyperlink $wsZh "F3" "https://github.com/OpenLocalizationTest/oltest/blob/915c598678366622ed8a59ae44ff321bc4968d68/e2e/75043d06-2288-40a2-a57e-05f54c1a449a.md" "75043d06-2288-40a2-a57e-05f54c1a449a.md"
Add-StyledHyperlink $wsZh "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ecd3281a4de4970c00f958a5185cbd454ee743f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/high/75043d06-2288-40a2-a57e-05f54c1a449a.bea90cac04fd87c8da04230db01b4947333f1089.zh-cn.xlf" "75043d06-2288-40a2-a57e-05f54c1a449a.bea90cac04fd87c8da04230db01b4947333f1089.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("H2").Value = "2016-03-30 11:28:32"
$wsDe.Range("H3").Value = "2016-03-30 11:28:32"

$hlA3de = Get-HyperlinkByAddr $wsDe '$A$3'
$a3AddressDe = $hlA3de.Address
$a3DisplayDe = $hlA3de.TextToDisplay
$hlA3de.Delete()

$hlD3de = Get-HyperlinkByAddr $wsDe '$D$3'
$d3AddressDe = $hlD3de.Address
$d3DisplayDe = $hlD3de.TextToDisplay
$hlD3de.Delete()

Add-StyledHyperlink $wsDe "F2" "https://github.com/OpenLocalizationTest/oltest/blob/915c598678366622ed8a59ae44ff321bc4968d68/e2e/27908121-ab1a-42e7-a0fa-0f05bdfe0441.md" "27908121-ab1a-42e7-a0fa-0f05bdfe0441.md"
Add-StyledHyperlink $wsDe "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a7c862e697f5cc8f00c4e36a628fa8dff5772b9/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/high/27908121-ab1a-42e7-a0fa-0f05bdfe0441.7677fbf96608daf076d6d3278a6411e68340f762.de-de.xlf" "27908121-ab1a-42e7-a0fa-0f05bdfe0441.7677fbf96608daf076d6d3278a6411e68340f762.de-de.xlf"

Add-StyledHyperlink $wsDe "A3" $a3AddressDe $a3DisplayDe
Add-StyledHyperlink $wsDe "D3" $d3AddressDe $d3DisplayDe

Add-StyledHyperlink $wsDe "F3" "https://github.com/OpenLocalizationTest/oltest/blob/915c598678366622ed8a59ae44ff321bc4968d68/e2e/75043d06-2288-40a2-a57e-05f54c1a449a.md" "75043d06-2288-40a2-a57e-05f54c1a449a.md"
Add-StyledHyperlink $wsDe "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a7c862e697f5cc8f00c4e36a628fa8dff5772b9/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/high/75043d06-2288-40a2-a57e-05f54c1a449a.bea90cac04fd87c8da04230db01b4947333f1089.de-de.xlf" "75043d06-2288-40a2-a57e-05f54c1a449a.bea90cac04fd87c8da04230db01b4947333f1089.de-de.xlf"

Write-Host "Handback report generated."
